# Auto-generated edit script updating the cryptos price table (columns B-E)
# to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.596.46"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.278.17"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "503.95"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.295.12"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0963"
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.341"
$ws.Range("E12").Value = "  +2.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.90"
$ws.Range("E13").Value = "  +4.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.28"
$ws.Range("E14").Value = "  +3.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.683.40"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "54.660.68"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.300.01"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.32"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "307.03"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.31"
$ws.Range("E24").Value = "  -2.95%  "
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.44"
$ws.Range("E27").Value = "  +1.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "170.64"
$ws.Range("E28").Value = "  -1.82%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0702"
$ws.Range("E29").Value = "  +2.05%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.05"
$ws.Range("E30").Value = "  +1.24%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.62"
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.12"
$ws.Range("E32").Value = "  +2.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.92"
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.910"
$ws.Range("E36").Value = "  -1.74%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.41"
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.375"
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.42"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.39"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "126.51"
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.82"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "250.02"
$ws.Range("E45").Value = "  +4.30%  "
$ws.Range("E46").Value = "  +0.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0901"
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.81"
$ws.Range("E51").Value = "  +0.38%  "
